$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.996.18"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "2.294.66"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Value = "'99.24"
$ws.Range("E6").Value = "  +1.88%  "
$ws.Range("D7").Value = "'0.504"
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("E9").Value = "  +2.31%  "
$ws.Range("E10").Value = "  +7.72%  "
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").Value = "'0.116"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").Value = "'18.13"
$ws.Range("E13").Value = "  +7.57%  "
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("D15").Value = "2.651.64"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").Value = "2.351.19"
$ws.Range("E16").Value = "  +2.73%  "
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("D18").Value = "42.893.78"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("E19").Value = "  +8.45%  "
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "'6.10"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").Value = "'67.76"
$ws.Range("D23").Value = "'236.01"
$ws.Range("E23").Value = "  -0.10%  "
$ws.Range("E24").Value = "  +9.71%  "
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").Value = "'24.92"
$ws.Range("E27").Value = "  +2.30%  "
$ws.Range("D29").Value = "'34.51"
$ws.Range("E29").Value = "  +1.80%  "
$ws.Range("D30").Value = "'167.30"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").Value = "'9.12"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "'5.01"
$ws.Range("E33").Value = "  +1.52%  "
$ws.Range("D34").Value = "'17.64"
$ws.Range("E34").Value = "  +3.71%  "
$ws.Range("D35").Value = "'4.62"
$ws.Range("E35").Value = "  -1.38%  "
$ws.Range("E36").Value = "  +1.18%  "
$ws.Range("D37").Value = "'0.0689"
$ws.Range("E37").Value = "  -0.71%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'1.79"
$ws.Range("E38").Value = "  +1.68%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.101"
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").Value = "'0.0291"
$ws.Range("E42").Value = "  +3.61%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").Value = "'2.29"
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("D44").Value = "1.972.61"
$ws.Range("E44").Value = "  -0.94%  "
$ws.Range("E45").Value = "  +3.29%  "
$ws.Range("D46").Value = "'2.90"
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("D47").Value = "'17.53"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").Value = "'55.47"
$ws.Range("E48").Value = "  +4.17%  "
$ws.Range("E49").Value = "  +3.52%  "
$ws.Range("D50").Value = "2.518.42"
$ws.Range("E50").Value = "  -0.08%  "
$ws.Range("D51").Value = "'70.75"
$ws.Range("E51").Value = "  +0.98%  "
